$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.940.68'
$ws.Range('E2').Value = '  -4.03%  '
$ws.Range('D3').Value = '1.640.78'
$ws.Range('E3').Value = '  -5.79%  '
$ws.Range('D4').Value = '''0.9992'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''232.46'
$ws.Range('E5').Value = '  -6.09%  '
$ws.Range('D6').Value = '''1.000'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').Value = '''0.4743'
$ws.Range('E7').Value = '  -5.69%  '
$ws.Range('D8').Value = '''0.2575'
$ws.Range('E8').Value = '  -6.16%  '
$ws.Range('D9').Value = '''0.06085'
$ws.Range('E9').Value = '  -1.74%  '
$ws.Range('D10').Value = '''0.07042'
$ws.Range('E10').Value = '  -3.20%  '
$ws.Range('D11').Value = '1.652.57'
$ws.Range('E11').Value = '  -5.14%  '
$ws.Range('D12').Value = '''14.51'
$ws.Range('E12').Value = '  -4.41%  '
$ws.Range('D13').Value = '''0.5880'
$ws.Range('E13').Value = '  -10.22%  '
$ws.Range('D14').Value = '''4.330'
$ws.Range('E14').Value = '  -7.77%  '
$ws.Range('D15').Value = '''73.59'
$ws.Range('E15').Value = '  -5.15%  '
$ws.Range('D16').Value = '''1.000'
$ws.Range('E16').Value = '  -0.01%  '
$ws.Range('D17').Value = '''1.001'
$ws.Range('E17').Value = '  +0.16%  '
$ws.Range('D18').Value = '24.939.72'
$ws.Range('E18').Value = '  -4.13%  '
$ws.Range('D19').Value = '''0.000006568'
$ws.Range('E19').Value = '  -4.14%  '
$ws.Range('D20').Value = '''11.21'
$ws.Range('E20').Value = '  -5.71%  '
$ws.Range('D21').Value = '1.857.01'
$ws.Range('E21').Value = '  -5.69%  '
$ws.Range('D22').Value = '''4.295'
$ws.Range('E22').Value = '  -5.94%  '
$ws.Range('D23').Value = '''8.528'
$ws.Range('E23').Value = '  -2.37%  '
$ws.Range('D24').Value = '''5.228'
$ws.Range('E24').Value = '  -3.11%  '
$ws.Range('D25').Value = '''133.59'
$ws.Range('E25').Value = '  -1.30%  '
$ws.Range('D26').Value = '''14.90'
$ws.Range('E26').Value = '  -2.35%  '
$ws.Range('D27').Value = '''1.379'
$ws.Range('E27').Value = '  -8.83%  '
$ws.Range('D28').Value = '''103.90'
$ws.Range('E28').Value = '  -1.36%  '
$ws.Range('D29').Value = '''1.635'
$ws.Range('E29').Value = '  -8.44%  '
$ws.Range('D30').Value = '''3.882'
$ws.Range('E30').Value = '  -1.84%  '
$ws.Range('D31').Value = '''0.07583'
$ws.Range('E31').Value = '  -6.95%  '
$ws.Range('D32').Value = '''3.559'
$ws.Range('E32').Value = '  -3.96%  '
$ws.Range('D33').Value = '''0.9996'
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('D34').Value = '''0.04262'
$ws.Range('E34').Value = '  -9.85%  '
$ws.Range('D35').Value = '''2.573'
$ws.Range('E35').Value = '  -3.54%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '''0.5910'
$ws.Range('E36').Value = '  -3.04%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = '''0.9241'
$ws.Range('E37').Value = '  -7.44%  '
$ws.Range('D38').Value = '''2.556'
$ws.Range('E38').Value = '  -7.25%  '
$ws.Range('D39').Value = '''0.8579'
$ws.Range('E39').Value = '  +6.24%  '
$ws.Range('D40').Value = '''0.9999'
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('B41').Value = 'PaxosStandard'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$ws.Range('D41').Value = '''1.000'
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = '''0.01494'
$ws.Range('E42').Value = '  -7.90%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = '''98.67'
$ws.Range('E43').Value = '  -2.54%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '''1.759'
$ws.Range('E44').Value = '  -9.10%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').Value = '''0.3693'
$ws.Range('E45').Value = '  -5.67%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '''4.666'
$ws.Range('E46').Value = '  -7.22%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '''0.1102'
$ws.Range('E47').Value = '  -5.82%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').Value = '''6.091'
$ws.Range('E48').Value = '  -4.28%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '''0.05199'
$ws.Range('E49').Value = '  -1.75%  '
$ws.Range('B50').Value = 'TrueUSD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range('D50').Value = '''1.002'
$ws.Range('E50').Value = '  -0.04%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').Value = '''28.69'
$ws.Range('E51').Value = '  -7.02%  '
